# "Alle Elemente für Animation 4b und 4c positioniert"
#
# Slide 7: duplicate the existing arrow connector (id=3, "Gerade Verbindung
# mit Pfeil 2") into a new connector ("Gerade Verbindung mit Pfeil 1"),
# reposition both, and move the original connector to the end of the
# z-order (after the new duplicate). Slide 8: reposition/resize three
# existing connectors.

$EMU_PER_PT = 12700.0

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 7 ("4b")
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# The original arrow connector (cNvPr id="3") is the first shape in the
# spTree. Duplicate it to create the new "Pfeil 1" connector.
$orig = $s7.Shapes.Item(1)

$dup = $orig.Duplicate()
$dup.Name = "Gerade Verbindung mit Pfeil 1"
$dup.HorizontalFlip = -1
$dup.Left = 3150082 / $EMU_PER_PT
$dup.Top = 1318437 / $EMU_PER_PT
$dup.Width = 4752149 / $EMU_PER_PT
$dup.Height = 1212112 / $EMU_PER_PT

# Move the original connector to the back of the creation order / front
# of the z-order stack so it ends up last (after the new duplicate).
$orig.ZOrder(0)   # msoBringToFront
$orig.Left = 3657989 / $EMU_PER_PT
$orig.Top = 2982202 / $EMU_PER_PT
$orig.Width = 0 / $EMU_PER_PT
$orig.Height = 2232837 / $EMU_PER_PT
$orig.Line.Weight = 46990 / $EMU_PER_PT

# ---------------------------------------------------------------------
# Slide 8 ("4c")
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

$sh2 = $s8.Shapes.Item(1)   # "Gerader Verbinder 1"
$sh2.Left = 6115050 / $EMU_PER_PT
$sh2.Top = 4554252 / $EMU_PER_PT
$sh2.Width = 5856412 / $EMU_PER_PT
$sh2.Height = 0 / $EMU_PER_PT

$sh3 = $s8.Shapes.Item(2)   # "Gerader Verbinder 2"
$sh3.Left = 6096000 / $EMU_PER_PT
$sh3.Top = 2177143 / $EMU_PER_PT
$sh3.Width = 0 / $EMU_PER_PT
$sh3.Height = 2397238 / $EMU_PER_PT

$sh4 = $s8.Shapes.Item(3)   # "Gerade Verbindung mit Pfeil 3"
$sh4.VerticalFlip = -1
$sh4.Left = 163043 / $EMU_PER_PT
$sh4.Top = 2157990 / $EMU_PER_PT
$sh4.Width = 5952007 / $EMU_PER_PT
$sh4.Height = 19153 / $EMU_PER_PT
